$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 307.22223
$ws.Range("I2").Value = 297
$ws.Range("K2").Value = 297
$ws.Range("M2").Value = -184
$ws.Range("H33").Value = 5263371
$ws.Range("I33").Value = 5555776.5
$ws.Range("K33").Value = 5555776.5
$ws.Range("M33").Value = -5555547.5
$ws.Range("H69").Value = 10286.5
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 10286.5
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("L69").Value = 30859.5
$ws.Range("N69").Value = -32607.5
$ws.Range("H72").Value = 10286.5
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 10286.5
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("L72").Value = 92578.5
$ws.Range("N72").Value = -101314.5
$ws.Range("H112").Value = 1416.5
$ws.Range("J112").Value = 2000
$ws.Range("L112").Value = 6000
$ws.Range("N112").Value = -8216
$ws.Range("H113").Value = 42128.234
$ws.Range("I113").Value = 78553.42999999999
$ws.Range("J113").Value = 10256.1875
$ws.Range("K113").Value = 78553.42999999999
$ws.Range("L113").Value = 10256.1875
$ws.Range("M113").Value = -75299.42999999999
$ws.Range("N113").Value = -16764.1875
$ws.Range("H132").Value = 9582101
$ws.Range("I132").Value = 10132504
$ws.Range("K132").Value = 30397512
$ws.Range("M132").Value = -30394982
$ws.Range("H138").Value = 3704.1428
$ws.Range("I138").Value = 1622.0588
$ws.Range("J138").Value = 4473.609
$ws.Range("K138").Value = 4866.1764
$ws.Range("L138").Value = 13420.827
$ws.Range("M138").Value = 273.8235999999997
$ws.Range("N138").Value = -23700.827

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 5015
$ws.Range("I39").Value = 5015
$ws.Range("K39").Value = 5015
$ws.Range("M39").Value = -4495
$ws.Range("H45").Value = 11444.379
$ws.Range("I45").Value = 9230.044
$ws.Range("J45").Value = 19932.666
$ws.Range("K45").Value = 9230.044
$ws.Range("L45").Value = 19932.666
$ws.Range("M45").Value = -8853.044
$ws.Range("N45").Value = -20686.666
$ws.Range("H74").Value = 62385.277
$ws.Range("I74").Value = 72681.28999999999
$ws.Range("K74").Value = 72681.28999999999
$ws.Range("M74").Value = -71807.28999999999
$ws.Range("H77").Value = 62385.277
$ws.Range("I77").Value = 72681.28999999999
$ws.Range("K77").Value = 363406.45
$ws.Range("M77").Value = -359038.45
$ws.Range("H122").Value = 1302.9333
$ws.Range("I122").Value = 1148.5
$ws.Range("J122").Value = 2306.75
$ws.Range("K122").Value = 3445.5
$ws.Range("L122").Value = 6920.25
$ws.Range("M122").Value = -995.5
$ws.Range("N122").Value = -11820.25
$ws.Range("H132").Value = 3301.7144
$ws.Range("I132").Value = 3185.3333
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 9555.999899999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -7025.999899999999
$ws.Range("N132").Value = -17060

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 13650
$ws.Range("J23").Value = 13650
$ws.Range("N23").Value = -14216
$ws.Range("L23").Value = 13650
$ws.Range("H64").Value = 893.7
$ws.Range("J64").Value = 530.6667
$ws.Range("L64").Value = 530.6667
$ws.Range("N64").Value = -980.6667
$ws.Range("H67").Value = 893.7
$ws.Range("J67").Value = 530.6667
$ws.Range("L67").Value = 530.6667
$ws.Range("N67").Value = -2090.6667
$ws.Range("H99").Value = 5865.6665
$ws.Range("I99").Value = 2043.2222
$ws.Range("K99").Value = 2043.2222
$ws.Range("M99").Value = -545.2221999999999
$ws.Range("H110").Value = 37500
$ws.Range("J110").Value = 37500
$ws.Range("N110").Value = -45680
$ws.Range("L110").Value = 37500

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1098.4615
$ws.Range("I16").Value = 920
$ws.Range("K16").Value = 920
$ws.Range("M16").Value = -633
$ws.Range("H28").Value = 38000
$ws.Range("J28").Value = 38000
$ws.Range("L28").Value = 38000
$ws.Range("N28").Value = -38490
$ws.Range("H31").Value = 33649.668
$ws.Range("I31").Value = 37326.395
$ws.Range("K31").Value = 37326.395
$ws.Range("M31").Value = -37031.395
$ws.Range("H34").Value = 33649.668
$ws.Range("I34").Value = 37326.395
$ws.Range("K34").Value = 37326.395
$ws.Range("M34").Value = -37124.395
$ws.Range("H58").Value = 2062.15
$ws.Range("I58").Value = 2173.2942
$ws.Range("K58").Value = 2173.2942
$ws.Range("M58").Value = -1970.2942
$ws.Range("H86").Value = 5533.3335
$ws.Range("I86").Value = 5533.3335
$ws.Range("K86").Value = 5533.3335
$ws.Range("M86").Value = -4410.3335
$ws.Range("H89").Value = 5533.3335
$ws.Range("I89").Value = 5533.3335
$ws.Range("K89").Value = 27666.6675
$ws.Range("M89").Value = -22050.6675
$ws.Range("H99").Value = 3944.4119
$ws.Range("I99").Value = 3805.7144
$ws.Range("K99").Value = 3805.7144
$ws.Range("M99").Value = -2307.7144
$ws.Range("H107").Value = 1718.5
$ws.Range("J107").Value = 1100
$ws.Range("L107").Value = 1100
$ws.Range("N107").Value = -4940
$ws.Range("H113").Value = 1098.4615
$ws.Range("I113").Value = 920
$ws.Range("K113").Value = 920
$ws.Range("M113").Value = 1250
$ws.Range("H126").Value = 3944.4119
$ws.Range("I126").Value = 3805.7144
$ws.Range("K126").Value = 11417.1432
$ws.Range("M126").Value = -8947.143199999999
$ws.Range("H136").Value = 2062.15
$ws.Range("I136").Value = 2173.2942
$ws.Range("K136").Value = 6519.882599999999
$ws.Range("M136").Value = -3969.882599999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3266162.5
$ws.Range("I4").Value = 3539404.5
$ws.Range("K4").Value = 10618213.5
$ws.Range("M4").Value = -10618101.5
$ws.Range("H114").Value = 2200
$ws.Range("J114").Value = 2180
$ws.Range("L114").Value = 6540
$ws.Range("N114").Value = -13048
$ws.Range("H131").Value = 46461.86
$ws.Range("I131").Value = 154624.53
$ws.Range("J131").Value = 8458.757
$ws.Range("K131").Value = 463873.59
$ws.Range("L131").Value = 25376.271
$ws.Range("M131").Value = -458833.59
$ws.Range("N131").Value = -35456.271

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1843.2
$ws.Range("I97").Value = 2118.1428
$ws.Range("J97").Value = 1201.6666
$ws.Range("K97").Value = 2118.1428
$ws.Range("L97").Value = 1201.6666
$ws.Range("M97").Value = -1622.1428
$ws.Range("N97").Value = -2193.6666
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("L105").Value = 0
$ws.Range("H113").Value = 7938.778
$ws.Range("I113").Value = 6983
$ws.Range("J113").Value = 8416.666999999999
$ws.Range("K113").Value = 6983
$ws.Range("L113").Value = 8416.666999999999
$ws.Range("M113").Value = -4813
$ws.Range("N113").Value = -12756.667
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("L135").Value = 0

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 1000000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H42").Value = 18505.5
$ws.Range("I42").Value = 25000
$ws.Range("J42").Value = 16340.667
$ws.Range("K42").Value = 25000
$ws.Range("L42").Value = 16340.667
$ws.Range("M42").Value = -24437
$ws.Range("N42").Value = -17466.667
$ws.Range("H49").Value = 18505.5
$ws.Range("I49").Value = 25000
$ws.Range("J49").Value = 16340.667
$ws.Range("K49").Value = 25000
$ws.Range("L49").Value = 16340.667
$ws.Range("M49").Value = -24853
$ws.Range("N49").Value = -16634.667
$ws.Range("H122").Value = 405556.44
$ws.Range("I122").Value = 504995.56
$ws.Range("K122").Value = 1514986.68
$ws.Range("M122").Value = -1512536.68
$ws.Range("H136").Value = 8799.9
$ws.Range("I136").Value = 8499.875
$ws.Range("K136").Value = 25499.625
$ws.Range("M136").Value = -22949.625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13912.571
$ws.Range("I45").Value = 11999
$ws.Range("J45").Value = 14231.5
$ws.Range("K45").Value = 11999
$ws.Range("L45").Value = 14231.5
$ws.Range("M45").Value = -11508
$ws.Range("N45").Value = -15213.5
$ws.Range("H100").Value = 1262
$ws.Range("I100").Value = 1126.9333
$ws.Range("J100").Value = 1464.6
$ws.Range("K100").Value = 2253.8666
$ws.Range("L100").Value = 2929.2
$ws.Range("M100").Value = -1712.8666
$ws.Range("N100").Value = -4011.2
$ws.Range("H110").Value = 150000
$ws.Range("J110").Value = 150000
$ws.Range("L110").Value = 150000
$ws.Range("N110").Value = -158180
$ws.Range("H122").Value = 3636
$ws.Range("I122").Value = 3206.182
$ws.Range("K122").Value = 9618.545999999998
$ws.Range("M122").Value = -7168.545999999998
$ws.Range("H136").Value = 2959.8462
$ws.Range("I136").Value = 2810.5
$ws.Range("J136").Value = 3198.8
$ws.Range("K136").Value = 8431.5
$ws.Range("L136").Value = 9596.400000000001
$ws.Range("M136").Value = -5881.5
$ws.Range("N136").Value = -14696.4
